$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("MPSP")
$ws.Range("C2").Value = 2.222894840016944
$ws.Range("C3").Value = 2.081351166747502
$ws.Range("C4").Value = 1.962467225995198
$ws.Range("C5").Value = 1.844379580943132
$ws.Range("C6").Value = 1.763225476374728
$ws.Range("F6").Value = 1.183398211708029
$ws.Range("C7").Value = 1.743076943634948
$ws.Range("F7").Value = 1.176159304169029
$ws.Range("C8").Value = 2.573321818083491
$ws.Range("C9").Value = 2.360463260563568
$ws.Range("C10").Value = 2.181077472300464
$ws.Range("C11").Value = 2.065482651075842
$ws.Range("C12").Value = 1.982096670058805
$ws.Range("F12").Value = 1.347707040680766
$ws.Range("C13").Value = 1.965543273948365
$ws.Range("F13").Value = 1.340645131982352
$ws.Range("C14").Value = 3.798504418416145
$ws.Range("C15").Value = 3.101354795070352
$ws.Range("C16").Value = 2.665321358710611
$ws.Range("C17").Value = 2.364075491314154
$ws.Range("C18").Value = 2.215306683784179
$ws.Range("F18").Value = 1.524438719629138
$ws.Range("C19").Value = 2.181626032584568
$ws.Range("F19").Value = 1.512136622448377

$ws = $wb.Worksheets.Item("GWP")
$ws.Range("C2").Value = -0.8633211588836323
$ws.Range("C3").Value = -1.17021592108262
$ws.Range("C4").Value = -1.425787129582262
$ws.Range("C5").Value = -1.600183394265429
$ws.Range("C6").Value = -1.711878685561119
$ws.Range("F6").Value = 3.78092009393867
$ws.Range("C7").Value = -1.743504222895545
$ws.Range("F7").Value = 3.773116407712585
$ws.Range("C8").Value = 0.1341207350706173
$ws.Range("C9").Value = -0.1298203378663828
$ws.Range("C10").Value = -0.3424236868735403
$ws.Range("C11").Value = -0.4929153403973364
$ws.Range("C12").Value = -0.5985278638278291
$ws.Range("F12").Value = 4.214732559935982
$ws.Range("C13").Value = -0.6286930730025446
$ws.Range("F13").Value = 4.201840746793652
$ws.Range("C14").Value = 0.8372937020526948
$ws.Range("C15").Value = 0.603239401840206
$ws.Range("C16").Value = 0.3943798609918671
$ws.Range("C17").Value = 0.2400555208856006
$ws.Range("C18").Value = 0.1517257797554393
$ws.Range("F18").Value = 4.844275418395246
$ws.Range("C19").Value = 0.1266036772861589
$ws.Range("F19").Value = 4.832343476480743

$ws = $wb.Worksheets.Item("COD Price")
$ws.Range("C2").Value = 73.48797373235243
$ws.Range("C3").Value = 51.30541777295414
$ws.Range("C4").Value = 35.51373386792045
$ws.Range("C5").Value = 25.3309272631408
$ws.Range("C6").Value = 15.78157702493059
$ws.Range("F6").Value = -64.59234436812989
$ws.Range("C7").Value = 12.8269031597602
$ws.Range("F7").Value = -69.08436617928808
$ws.Range("C8").Value = 122.5309115816169
$ws.Range("C9").Value = 86.32869884435117
$ws.Range("C10").Value = 58.48920111832051
$ws.Range("C11").Value = 40.70701943842687
$ws.Range("C12").Value = 27.44440028938203
$ws.Range("F12").Value = -44.76062896011862
$ws.Range("C13").Value = 23.48474074193043
$ws.Range("F13").Value = -50.41848222624078
$ws.Range("C14").Value = 405.9016141644476
$ws.Range("C15").Value = 264.3052443264097
$ws.Range("C16").Value = 158.8977714287114
$ws.Range("C17").Value = 83.72441250124015
$ws.Range("C18").Value = 39.5382757341325
$ws.Range("F18").Value = -26.5366405259764
$ws.Range("C19").Value = 34.35855151237573
$ws.Range("F19").Value = -35.10303027300269

$ws = $wb.Worksheets.Item("COD GWP")
$ws.Range("C2").Value = -1174.01788651092
$ws.Range("C3").Value = -1183.243328483077
$ws.Range("C4").Value = -1191.026939031916
$ws.Range("C5").Value = -1187.234103477948
$ws.Range("C6").Value = -1184.77899704568
$ws.Range("F6").Value = -1678.572680181781
$ws.Range("C7").Value = -1183.806756452163
$ws.Range("F7").Value = -1687.079993241087
$ws.Range("C8").Value = -996.8199910720191
$ws.Range("C9").Value = -993.1550244986875
$ws.Range("C10").Value = -993.6949649005603
$ws.Range("C11").Value = -993.2773402229712
$ws.Range("C12").Value = -990.1938840079094
$ws.Range("F12").Value = -1413.969643110279
$ws.Range("C13").Value = -988.5667025547648
$ws.Range("F13").Value = -1421.018183103475
$ws.Range("C14").Value = -833.9761292759695
$ws.Range("C15").Value = -845.3732885747194
$ws.Range("C16").Value = -846.4284952975937
$ws.Range("C17").Value = -844.819853559672
$ws.Range("C18").Value = -846.7576401731577
$ws.Range("F18").Value = -1179.894970869601
$ws.Range("C19").Value = -846.5665255859478
$ws.Range("F19").Value = -1185.956082599678
